$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "30 January"
$ws.Range("B3").Value = "EXA club"

$ws.Range("C3").Select()
